# Femacal de La Calera - Poroto verde: insert a new weekly record.
# A new row is inserted at row 278 (shifting the existing rows 278-298
# down to 279-299) and populated with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 278, pushing rows 278:298 down to 279:299.
$ws.Rows("278:278").Insert()

# Populate the newly inserted row 278 with the new record's values.
$ws.Range("A278").Value = 3
$ws.Range("B278").Value = "Femacal de La Calera"
$ws.Range("C278").Value = "Coquimbo"
$ws.Range("D278").Value = 44585
$ws.Range("E278").Value = 5
$ws.Range("F278").Value = 100112031
$ws.Range("G278").Value = "Poroto verde"
$ws.Range("H278").Value = "Magnum"
$ws.Range("I278").Value = "Primera"
$ws.Range("J278").Value = 73
$ws.Range("K278").Value = 35000
$ws.Range("L278").Value = 36000
$ws.Range("M278").Value = 35479
$ws.Range("N278").Value = '$/malla 25 kilos'
$ws.Range("O278").Value = "Provincia de Santiago"
$ws.Range("P278").Value = 1419
$ws.Range("Q278").Value = 25
$ws.Range("R278").Value = "Hortaliza"
